$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3183
$ws1.Range("F5").Value = 6867
$ws1.Range("F6").Value = 2002
$ws1.Range("F8").Value = 72
$ws1.Range("F10").Value = 35
$ws1.Range("F11").Value = 68
$ws1.Range("F12").Value = 21
$ws1.Range("F13").Value = 147
$ws1.Range("F14").Value = 179

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 16

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3183
$ws4.Range("F3").Value = 16
$ws4.Range("F6").Value = 6867
$ws4.Range("F7").Value = 2002
$ws4.Range("F9").Value = 72
$ws4.Range("F11").Value = 35
$ws4.Range("F12").Value = 68
$ws4.Range("F13").Value = 21
$ws4.Range("F14").Value = 147
$ws4.Range("F15").Value = 179

$wb.Save()
